$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title: merge "Super Hyper Rem " + "Lezar" (spell-checked) + " Turbo:"
#    into a single run "Super Hyper Rem Lezar Turbo:" (drops the proofErr
#    spell-check wrapper around "Lezar").
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Super Hyper Rem Lezar Turbo:", $false, $false, $false, $false, $false,
    $true, 1, $false, "Super Hyper Rem Lezar Turbo:", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "Mechanics:" -> "Mechanics" (drop trailing colon)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Mechanics:", $false, $false, $false, $false, $false,
    $true, 1, $false, "Mechanics", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Features section: re-split the runs for the "-Space Flight",
#    "-Different power ups" and "-Boss" bullets (visible text is unchanged).
# ---------------------------------------------------------------------------
$d2 = $word.ActiveDocument
for ($i = 1; $i -le $d2.Paragraphs.Count; $i++) {
    $para = $d2.Paragraphs($i)
    if ($para.Range.Text.Contains("Space Flight")) {
        $spaceFlightIdx = $i
    }
}
$pSpaceFlight = $d2.Paragraphs($spaceFlightIdx)
$pDifferent = $d2.Paragraphs($spaceFlightIdx + 1)
$pBoss = $d2.Paragraphs($spaceFlightIdx + 2)

$rngFeatures = $d2.Range($pSpaceFlight.Range.Start, $pBoss.Range.End)
$featuresXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/></w:r><w:r><w:t xml:space="preserve">-Space Flight </w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>-Different power ups</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t>-</w:t></w:r><w:r><w:t>Boss</w:t></w:r></w:p>'
$rngFeatures.InsertXML($featuresXml)

# ---------------------------------------------------------------------------
# 4) Remove the "Technical:" heading, the "Detail any numerical value..."
#    paragraph, and the blank paragraph that followed it (they moved out of
#    this document entirely).
# ---------------------------------------------------------------------------
$d3 = $word.ActiveDocument
$techStart = -1
$techEnd = -1
for ($i = 1; $i -le $d3.Paragraphs.Count; $i++) {
    $para = $d3.Paragraphs($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "Technical:") {
        $techStart = $i
    }
    if ($para.Range.Text.Contains("Detail any numerical value")) {
        $techEnd = $i + 1
    }
}
$pTechStart = $d3.Paragraphs($techStart)
$pTechEnd = $d3.Paragraphs($techEnd)
$rngTech = $d3.Range($pTechStart.Range.Start, $pTechEnd.Range.End)
$rngTech.Delete()

# ---------------------------------------------------------------------------
# 5) Replace the blank paragraph that used to sit between "-Power ups
#    (Boosted Attack)" and "Platform:" with the large new set of sections:
#    Enemy Types, Music, Progression, Level, Drive, Art Genre.
# ---------------------------------------------------------------------------
$d4 = $word.ActiveDocument
$powerUpsIdx = -1
for ($i = 1; $i -le $d4.Paragraphs.Count; $i++) {
    $para = $d4.Paragraphs($i)
    if ($para.Range.Text.Contains("Power ups (Boosted Attack)")) {
        $powerUpsIdx = $i
    }
}
$pBlank = $d4.Paragraphs($powerUpsIdx + 1)
$bigXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Enemy Types</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:t>-Scarab: The Scarab type ship will attempt to fly into you and if they succeed in doing so they will detonate on contact.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>-Killer Bee:</w:t></w:r><w:r><w:t xml:space="preserve"> The Killer Bee type ship fly will fly in wave like pattern across the screen and shoot burst from their guns when they face the middle of the screens</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:t xml:space="preserve">-Wasp: The Wasp type ship stays at </w:t></w:r><w:r><w:t>the top of the screen firing a beam that goes across the screen and then the ship strafes slightly to the left or right for a few seconds the</w:t></w:r><w:r><w:t>n stops and flies straight down</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Music</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr><w:r><w:tab/><w:t>-</w:t></w:r><w:r><w:t>Science Fantasy</w:t></w:r><w:r><w:t xml:space="preserve"> Electronica</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Progression</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r><w:t xml:space="preserve">-You progress through </w:t></w:r><w:r><w:t>the levels by defeating the enemies until you reach</w:t></w:r><w:r><w:t xml:space="preserve"> the boss of the level. Once the boss is defeated you progress on to the next level. You collect power ups throughout the game to make your ships and weapons stronger. Once you go through all of the levels and defeat the final boss the game ends.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Level</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r><w:t>-</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>As you progress through each level the game gets progressively ha</w:t></w:r><w:r><w:t>rder, until you reach the end and beat the game</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Drive</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r><w:t xml:space="preserve">-The drive for the game </w:t></w:r><w:r><w:t>is to get stronger and beat the more challenging levels until you finally reach the end of the game</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Art Genre</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr><w:r><w:tab/><w:t>-</w:t></w:r><w:r><w:t>Realistic</w:t></w:r></w:p>'
$pBlank.Range.InsertXML($bigXml)

# ---------------------------------------------------------------------------
# 6) Remove the old "_GoBack" bookmark that used to sit right after
#    "Alien Bugs: Bad" (it now lives in the new "Drive" paragraph instead).
#    Also drop the stray <w:lastRenderedPageBreak/> that used to precede
#    "Target Audience:" (it moved to the new "-Killer Bee:" run above).
# ---------------------------------------------------------------------------
$d5 = $word.ActiveDocument
for ($i = 1; $i -le $d5.Paragraphs.Count; $i++) {
    $para = $d5.Paragraphs($i)
    if ($para.Range.Text.Contains("Alien Bugs")) {
        $bad = $para.Range.Text.TrimEnd([char]13, [char]7)
        $para.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Ali</w:t></w:r><w:r><w:t>en Bugs</w:t></w:r><w:r><w:t>: Bad</w:t></w:r></w:p>')
        break
    }
}

$d6 = $word.ActiveDocument
for ($i = 1; $i -le $d6.Paragraphs.Count; $i++) {
    $para = $d6.Paragraphs($i)
    if ($para.Range.Text.Contains("Target Audience")) {
        $para.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="480" w:lineRule="auto"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Target Audience: </w:t></w:r></w:p>')
        break
    }
}

# ---------------------------------------------------------------------------
# 7) Team Members / Jobs section: merge the Ryan / Barry / Chives runs.
# ---------------------------------------------------------------------------
$d7 = $word.ActiveDocument
for ($i = 1; $i -le $d7.Paragraphs.Count; $i++) {
    $para = $d7.Paragraphs($i)
    $t = $para.Range.Text
    if ($t.StartsWith("Ryan")) {
        $para.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Ryan: Ship Models/Textures</w:t></w:r><w:r><w:t>/Programming(Misc.)</w:t></w:r></w:p>')
    }
    elseif ($t.StartsWith("Barry")) {
        $para.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Barry: Environ</w:t></w:r><w:r><w:t>ment</w:t></w:r><w:r><w:t xml:space="preserve"> Models</w:t></w:r><w:r><w:t>/Textures/</w:t></w:r><w:r><w:t>Audio</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>/Programming(Misc.)</w:t></w:r></w:p>')
    }
    elseif ($t.StartsWith("Chives")) {
        $para.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Chives: Programming (</w:t></w:r><w:r><w:t>Main</w:t></w:r><w:r><w:t>)</w:t></w:r></w:p>')
    }
}

Write-Output "edit complete"
